$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")
$ws.Range("F7").Value = 124
$ws.Range("F9").Value = 1.076
